$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: source data from original row 18
$ws.Range("D2").Value = 44789
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = "$/bandeja 18 kilos"
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18

# Row 3: source data from original row 12
$ws.Range("D3").Value = 44616
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("S3").Value = 917

# Row 4: source data from original row 10
$ws.Range("D4").Value = 44784
$ws.Range("M4").Value = 300

# Row 5: source data from original row 11
$ws.Range("D5").Value = 44489
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 26000
$ws.Range("O5").Value = 27000
$ws.Range("P5").Value = 26500
$ws.Range("S5").Value = 1472

# Row 6: source data from original row 19
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 270
$ws.Range("Q6").Value = "$/bandeja 18 kilos"

# Row 7: source data from original row 8
$ws.Range("D7").Value = 44487
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("S7").Value = 1450

# Row 8: source data from original row 3
$ws.Range("D8").Value = 44614
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("Q8").Value = "$/bandeja 18 kilos"
$ws.Range("S8").Value = 1139
$ws.Range("T8").Value = 18

# Row 9: source data from original row 13
$ws.Range("D9").Value = 44602
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 270
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = "$/bandeja 18 kilos"
$ws.Range("S9").Value = 1139
$ws.Range("T9").Value = 18

# Row 10: source data from original row 9
$ws.Range("D10").Value = 44673
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 400
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("S10").Value = 1450
$ws.Range("T10").Value = 10

# Row 11: source data from original row 17
$ws.Range("D11").Value = 44418
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 10500
$ws.Range("Q11").Value = "$/bandeja 10 kilos"
$ws.Range("S11").Value = 1050
$ws.Range("T11").Value = 10

# Row 12: source data from original row 6
$ws.Range("D12").Value = 44263
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = "$/caja 18 kilos"
$ws.Range("S12").Value = 1194

# Row 13: source data from original row 4
$ws.Range("D13").Value = 44656
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("S13").Value = 1083

# Row 17: source data from original row 2
$ws.Range("D17").Value = 44706
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 9000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 9500
$ws.Range("S17").Value = 950

# Row 18: source data from original row 7
$ws.Range("D18").Value = 44819
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 17500
$ws.Range("Q18").Value = "$/bandeja 10 kilos"
$ws.Range("S18").Value = 1750
$ws.Range("T18").Value = 10

# Row 19: source data from original row 5
$ws.Range("D19").Value = 44291
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 17500
$ws.Range("S19").Value = 972

